$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.511.35'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '3.766.32'
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.68'
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.23'
$ws.Range("E6").Value = '  +2.22%  '

$ws.Range("D7").Value = '3.762.55'
$ws.Range("E7").Value = '  -0.90%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  -1.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("E12").Value = '  -3.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.25'
$ws.Range("E13").Value = '  -1.10%  '

$ws.Range("E14").Value = '  -0.77%  '

$ws.Range("D15").Value = '4.395.80'
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("D16").Value = '3.764.34'
$ws.Range("E16").Value = '  -0.99%  '

$ws.Range("D17").Value = '70.539.71'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("E18").Value = '  -2.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.56'
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.64'
$ws.Range("E20").Value = '  -0.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.29'
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.21'
$ws.Range("E22").Value = '  -3.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  -2.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.63'
$ws.Range("E24").Value = '  +5.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.29'
$ws.Range("E25").Value = '  -1.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.13'
$ws.Range("E26").Value = '  -3.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.39'
$ws.Range("E27").Value = '  +4.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000134'
$ws.Range("E28").Value = '  +5.68%  '

$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.49'
$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.94'
$ws.Range("E31").Value = '  +2.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.94'
$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.65'
$ws.Range("E33").Value = '  -4.04%  '

$ws.Range("E34").Value = '  -0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").Value = '  +1.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.15'
$ws.Range("E37").Value = '  -1.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.356'
$ws.Range("E38").Value = '  +3.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.142'
$ws.Range("E39").Value = '  +7.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.11'
$ws.Range("E40").Value = '  +14.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.09'
$ws.Range("E41").Value = '  -4.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '50.16'
$ws.Range("E42").Value = '  -2.72%  '

$ws.Range("B43").Value = 'Arweave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '45.67'
$ws.Range("E43").Value = '  +2.59%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '437.46'
$ws.Range("E44").Value = '  +2.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.64'
$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").Value = '2.960.09'
$ws.Range("E46").Value = '  -5.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0364'
$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.50'
$ws.Range("E48").Value = '  -1.50%  '

$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.84'
$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("E51").Value = '  -0.02%  '
